$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "( {`n`t""message"": ""TOKEN: abc, TOKEN: def, TOKEN: ghi, REGEX:[0-9a-z]{4}-[0-9a-z]{4}-[0-9a-z]{4}-[0-9a-z]{4}, VAR: var1, REGEX:Bar[0-9]{2}-sector[0-9]{2}, VAR: var2 ""`n} AND AFTER 00:00:00:10 {`n`t""message"": ""TOKEN:Node down, VAR:{var1}""`n} )"

$ws.Range("C4").Value = $newText

$ws.Range("C4").Select()
